$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper colors
$blackColor = 0          # RGB(0,0,0)
$redColor   = 1974729    # RGB(201,33,30) -> 0xC9211E

function Set-TwoRunHeader($addr, $plainText, $starText) {
    $cell = $ws.Range($addr)
    $full = $plainText + $starText
    $cell.Value2 = $full

    $plainLen = $plainText.Length
    $starLen  = $starText.Length

    $run1 = $cell.Characters(1, $plainLen)
    $run1.Font.Name = "Calibri"
    $run1.Font.Size = 11
    $run1.Font.Bold = $true
    $run1.Font.Color = $blackColor

    $run2 = $cell.Characters($plainLen + 1, $starLen)
    $run2.Font.Name = "Calibri"
    $run2.Font.Size = 11
    $run2.Font.Bold = $true
    $run2.Font.Color = $redColor
}

# Gender, Designation, Phone become two-run "required" headers (bold black text + bold red asterisk)
Set-TwoRunHeader "C1" "Gender " "*"
Set-TwoRunHeader "D1" "Designation " "*"
Set-TwoRunHeader "E1" "Phone " "*"

# Email / Address simply gain a trailing space, formatting unchanged
$ws.Range("F1").Value2 = "Email "
$ws.Range("G1").Value2 = "Address "

# Remove the stray row 2 (B2 held a lone space string) entirely
$ws.Range("B2").EntireRow.Delete()

# Update the active selection to C1 (was H1)
$ws.Range("C1").Select()
